# Auto-generated edit script: updates market-price-derived columns (H-N)
# on the Leve-profitability sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 271
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H62").Value = 2599.5
$ws.Range("I62").Value = 2685.1428
$ws.Range("K62").Value = 2685.1428
$ws.Range("M62").Value = -2061.1428
$ws.Range("H65").Value = 2599.5
$ws.Range("I65").Value = 2685.1428
$ws.Range("K65").Value = 13425.714
$ws.Range("M65").Value = -10305.714
$ws.Range("H106").Value = 6805.1177
$ws.Range("I106").Value = 5479.2
$ws.Range("K106").Value = 5479.2
$ws.Range("M106").Value = -4848.2
$ws.Range("H125").Value = 685
$ws.Range("J125").Value = 677.5
$ws.Range("L125").Value = 6097.5
$ws.Range("N125").Value = -11017.5
$ws.Range("H138").Value = 5000
$ws.Range("J138").Value = 5000
$ws.Range("L138").Value = 15000
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1081.1666
$ws.Range("J5").Value = 236.625
$ws.Range("L5").Value = 236.625
$ws.Range("N5").Value = -460.625
$ws.Range("H102").Value = 15874220
$ws.Range("I102").Value = 20409452
$ws.Range("J102").Value = 907
$ws.Range("K102").Value = 20409452
$ws.Range("L102").Value = 907
$ws.Range("M102").Value = -20407830
$ws.Range("N102").Value = -4151
$ws.Range("H110").Value = 7095.1665
$ws.Range("I110").Value = 7095.1665
$ws.Range("K110").Value = 7095.1665
$ws.Range("M110").Value = -5050.1665
$ws.Range("H122").Value = 2177.2727
$ws.Range("I122").Value = 2295
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6885
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -4435
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 3916.6667
$ws.Range("I132").Value = 1750
$ws.Range("K132").Value = 5250
$ws.Range("M132").Value = -2720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1081.1666
$ws.Range("J4").Value = 236.625
$ws.Range("L4").Value = 236.625
$ws.Range("N4").Value = -466.625
$ws.Range("H22").Value = 321.85715
$ws.Range("I22").Value = 321.85715
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 321.85715
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -148.85715
$ws.Range("H86").Value = 1750
$ws.Range("I86").Value = 1750
$ws.Range("K86").Value = 1750
$ws.Range("M86").Value = -627
$ws.Range("H89").Value = 1750
$ws.Range("I89").Value = 1750
$ws.Range("K89").Value = 8750
$ws.Range("M89").Value = -3134
$ws.Range("H105").Value = 90912936
$ws.Range("I105").Value = 800
$ws.Range("K105").Value = 800
$ws.Range("M105").Value = 947
$ws.Range("H107").Value = 1938.5714
$ws.Range("I107").Value = 1938.5714
$ws.Range("K107").Value = 1938.5714
$ws.Range("M107").Value = -18.57140000000004
$ws.Range("H134").Value = 2874.1667
$ws.Range("I134").Value = 2874.1667
$ws.Range("K134").Value = 8622.500100000001
$ws.Range("M134").Value = -6087.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1311.9
$ws.Range("I113").Value = 383.33334
$ws.Range("J113").Value = 1709.8572
$ws.Range("K113").Value = 1150.00002
$ws.Range("L113").Value = 5129.571599999999
$ws.Range("M113").Value = 1019.99998
$ws.Range("N113").Value = -9469.571599999999
$ws.Range("H139").Value = 4485.7144
$ws.Range("I139").Value = 6466.6665
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 19399.9995
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -14259.9995
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7525.4
$ws.Range("I70").Value = 6337.4
$ws.Range("J70").Value = 8119.4
$ws.Range("K70").Value = 6337.4
$ws.Range("L70").Value = 8119.4
$ws.Range("M70").Value = -6067.4
$ws.Range("N70").Value = -8659.4
$ws.Range("H73").Value = 7525.4
$ws.Range("I73").Value = 6337.4
$ws.Range("J73").Value = 8119.4
$ws.Range("K73").Value = 6337.4
$ws.Range("L73").Value = 8119.4
$ws.Range("M73").Value = -5401.4
$ws.Range("N73").Value = -9991.4
$ws.Range("H80").Value = 1838.1875
$ws.Range("I80").Value = 1785.4445
$ws.Range("J80").Value = 1906
$ws.Range("K80").Value = 1785.4445
$ws.Range("L80").Value = 1906
$ws.Range("M80").Value = -787.4445000000001
$ws.Range("N80").Value = -3902
$ws.Range("H83").Value = 1838.1875
$ws.Range("I83").Value = 1785.4445
$ws.Range("J83").Value = 1906
$ws.Range("K83").Value = 8927.2225
$ws.Range("L83").Value = 9530
$ws.Range("M83").Value = -3935.2225
$ws.Range("N83").Value = -19514
$ws.Range("H97").Value = 470.23077
$ws.Range("I97").Value = 380.7
$ws.Range("K97").Value = 380.7
$ws.Range("M97").Value = 115.3
$ws.Range("H102").Value = 3113.5454
$ws.Range("J102").Value = 1374.5
$ws.Range("L102").Value = 1374.5
$ws.Range("N102").Value = -4618.5
$ws.Range("H122").Value = 2406.889
$ws.Range("I122").Value = 1809.8422
$ws.Range("J122").Value = 3824.875
$ws.Range("K122").Value = 5429.5266
$ws.Range("L122").Value = 11474.625
$ws.Range("M122").Value = -2979.5266
$ws.Range("N122").Value = -16374.625
$ws.Range("H123").Value = 18155
$ws.Range("J123").Value = 18155
$ws.Range("L123").Value = 18155
$ws.Range("N123").Value = -23055
$ws.Range("H132").Value = 3816.3333
$ws.Range("I132").Value = 2966.3333
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 8898.999899999999
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -6368.999899999999
$ws.Range("N132").Value = -19059.0005
$ws.Range("H133").Value = 70780
$ws.Range("J133").Value = 70780
$ws.Range("L133").Value = 70780
$ws.Range("N133").Value = -80900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1845.5294
$ws.Range("I22").Value = 1481.6666
$ws.Range("J22").Value = 2718.8
$ws.Range("K22").Value = 1481.6666
$ws.Range("L22").Value = 2718.8
$ws.Range("M22").Value = -1186.6666
$ws.Range("N22").Value = -3308.8
$ws.Range("H27").Value = 1845.5294
$ws.Range("I27").Value = 1481.6666
$ws.Range("J27").Value = 2718.8
$ws.Range("K27").Value = 1481.6666
$ws.Range("L27").Value = 2718.8
$ws.Range("M27").Value = -1374.6666
$ws.Range("N27").Value = -2932.8
$ws.Range("H61").Value = 2322.9412
$ws.Range("I61").Value = 2040.6154
$ws.Range("K61").Value = 2040.6154
$ws.Range("M61").Value = -1838.6154
$ws.Range("H68").Value = 8473.666999999999
$ws.Range("I68").Value = 7022
$ws.Range("K68").Value = 7022
$ws.Range("M68").Value = -6273
$ws.Range("H71").Value = 8473.666999999999
$ws.Range("I71").Value = 7022
$ws.Range("K71").Value = 35110
$ws.Range("M71").Value = -31366
$ws.Range("H82").Value = 2462.125
$ws.Range("I82").Value = 2424.75
$ws.Range("J82").Value = 2499.5
$ws.Range("K82").Value = 2424.75
$ws.Range("L82").Value = 2499.5
$ws.Range("M82").Value = -2063.75
$ws.Range("N82").Value = -3221.5
$ws.Range("H85").Value = 2462.125
$ws.Range("I85").Value = 2424.75
$ws.Range("J85").Value = 2499.5
$ws.Range("K85").Value = 2424.75
$ws.Range("L85").Value = 2499.5
$ws.Range("M85").Value = -1176.75
$ws.Range("N85").Value = -4995.5
$ws.Range("H113").Value = 2322.9412
$ws.Range("I113").Value = 2040.6154
$ws.Range("K113").Value = 2040.6154
$ws.Range("M113").Value = 129.3846000000001
$ws.Range("H128").Value = 84950
$ws.Range("J128").Value = 84950
$ws.Range("L128").Value = 84950
$ws.Range("N128").Value = -94910
$ws.Range("H132").Value = 3020.611
$ws.Range("I132").Value = 2871.6
$ws.Range("K132").Value = 8614.799999999999
$ws.Range("M132").Value = -6084.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 23999
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H122").Value = 1444.1428
$ws.Range("I122").Value = 1434.8334
$ws.Range("K122").Value = 4304.5002
$ws.Range("M122").Value = -1854.5002
$ws.Range("H132").Value = 2417.1924
$ws.Range("I132").Value = 1960.6316
$ws.Range("K132").Value = 5881.8948
$ws.Range("M132").Value = -3351.8948
